$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Fecha(serial), Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Origen, PrecioKg
$rowsData = @(
    ,@(2, 44161, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(3, 44876, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(4, 44230, "Primera", 16000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(5, 44880, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(6, 44167, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(7, 44875, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(8, 44881, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(9, 44874, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(10, 44181, "Primera", 12000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(11, 44882, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(12, 44847, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(13, 44210, "Primera", 8800, 2500, 3000, 2750, "Provincia de Chacabuco", 28)
    ,@(14, 44232, "Primera", 16000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(15, 44902, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(16, 44873, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(17, 44229, "Primera", 16000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(18, 44188, "Primera", 12000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(19, 44859, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(20, 44245, "Primera", 9000, 3000, 3000, 3000, "Región Metropolitana", 30)
    ,@(21, 44245, "Segunda", 5000, 2500, 2500, 2500, "Región Metropolitana", 25)
    ,@(22, 44204, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(23, 44162, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(24, 44855, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(25, 44602, "Primera", 12000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(26, 44602, "Segunda", 6000, 2500, 2500, 2500, "Provincia de Chacabuco", 25)
    ,@(27, 44600, "Primera", 1300, 3500, 4000, 3808, "Región Metropolitana", 38)
    ,@(28, 44168, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(29, 44215, "Primera", 16000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(30, 44186, "Primera", 10000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(31, 44160, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(32, 44883, "Primera", 9700, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(33, 44901, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(34, 44214, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(35, 44189, "Primera", 16000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(36, 44209, "Primera", 7000, 2500, 3000, 2750, "Provincia de Chacabuco", 28)
    ,@(37, 44187, "Primera", 12000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(38, 44159, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(39, 44166, "Primera", 7000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(40, 44231, "Primera", 12000, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(41, 44860, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(42, 44845, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
    ,@(43, 44846, "Primera", 7900, 3000, 3000, 3000, "Provincia de Chacabuco", 30)
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Range("D$r").Value2 = $row[1]
    $ws.Range("I$r").Value2 = $row[2]
    $ws.Range("J$r").Value2 = $row[3]
    $ws.Range("K$r").Value2 = $row[4]
    $ws.Range("L$r").Value2 = $row[5]
    $ws.Range("M$r").Value2 = $row[6]
    $ws.Range("O$r").Value2 = $row[7]
    $ws.Range("P$r").Value2 = $row[8]
}

Write-Host "done"